$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.060.22"
$ws.Range("E2").Value = "  -2.27%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.798.06"
$ws.Range("E3").Value = "  -2.86%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.43%  "

# Row 5 - now USDC (was BNB)
$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").Value = "'1.003"
$ws.Range("E5").Value = "  +0.35%  "

# Row 6 - now BNB (was USDC)
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'307.45"
$ws.Range("E6").Value = "  -2.26%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4175"
$ws.Range("E7").Value = "  -2.45%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3544"
$ws.Range("E8").Value = "  -4.19%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.07064"
$ws.Range("E9").Value = "  -4.07%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "'0.8421"
$ws.Range("E10").Value = "  -3.83%  "

# Row 11 - Solana
$ws.Range("D11").Value = "'20.18"
$ws.Range("E11").Value = "  -3.83%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.862.96"
$ws.Range("E12").Value = "  +1.42%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'5.273"
$ws.Range("E13").Value = "  -3.38%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "'6.339"
$ws.Range("E14").Value = "  -3.69%  "

# Row 15 - TRON
$ws.Range("D15").Value = "'0.06756"
$ws.Range("E15").Value = "  -2.99%  "

# Row 16 - BinanceUSD
$ws.Range("D16").Value = "'1.008"
$ws.Range("E16").Value = "  +0.59%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "'79.87"
$ws.Range("E17").Value = "  -1.22%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "'0.000008691"
$ws.Range("E18").Value = "  -4.48%  "

# Row 19 - Dai
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.17%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "'15.03"
$ws.Range("E20").Value = "  -3.47%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "27.012.03"
$ws.Range("E21").Value = "  -2.27%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'5.046"
$ws.Range("E22").Value = "  -0.90%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "'10.93"
$ws.Range("E23").Value = "  -0.84%  "

# Row 24 - WrappedliquidstakedEther2.0
$ws.Range("D24").Value = "2.003.81"
$ws.Range("E24").Value = "  -2.43%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "'1.947"
$ws.Range("E25").Value = "  -0.70%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'153.14"
$ws.Range("E26").Value = "  -1.02%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'18.08"
$ws.Range("E27").Value = "  -2.59%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "'4.988"
$ws.Range("E28").Value = "  -6.03%  "

# Row 29 - BitcoinCash
$ws.Range("D29").Value = "'112.94"
$ws.Range("E29").Value = "  -2.03%  "

# Row 30 - LidoDAOToken
$ws.Range("D30").Value = "'1.642"
$ws.Range("E30").Value = "  -12.07%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "'0.08899"
$ws.Range("E31").Value = "  -0.09%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "'0.7168"
$ws.Range("E32").Value = "  -8.66%  "

# Row 33 - HuobiToken
$ws.Range("D33").Value = "'2.858"
$ws.Range("E33").Value = "  -4.04%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "'4.320"
$ws.Range("E34").Value = "  -6.07%  "

# Row 35 - Frax
$ws.Range("D35").Value = "'1.004"
$ws.Range("E35").Value = "  +0.38%  "

# Row 36 - ARBITRUM
$ws.Range("D36").Value = "'1.074"
$ws.Range("E36").Value = "  -7.34%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").Value = "'1.075"
$ws.Range("E37").Value = "  -3.00%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "'0.01892"
$ws.Range("E38").Value = "  -3.27%  "

# Row 39 - Hedera
$ws.Range("D39").Value = "'0.05101"
$ws.Range("E39").Value = "  -5.91%  "

# Row 40 - Algorand
$ws.Range("D40").Value = "'0.1615"
$ws.Range("E40").Value = "  -3.65%  "

# Row 41 - TheSandbox
$ws.Range("D41").Value = "'0.4925"
$ws.Range("E41").Value = "  -4.89%  "

# Row 42 - MXToken
$ws.Range("D42").Value = "'2.586"
$ws.Range("E42").Value = "  -8.88%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "'6.087"
$ws.Range("E43").Value = "  -9.82%  "

# Row 44 - Aptos
$ws.Range("D44").Value = "'8.022"
$ws.Range("E44").Value = "  -7.38%  "

# Row 45 - now Quant (was EnergySwap)
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'104.30"
$ws.Range("E45").Value = "  -2.86%  "

# Row 46 - PaxDollar
$ws.Range("D46").Value = "'1.003"
$ws.Range("E46").Value = "  +0.32%  "

# Row 47 - now EnergySwap (was Quant)
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.17"
$ws.Range("E47").Value = "  -3.75%  "

# Row 48 - Cronos
$ws.Range("D48").Value = "'0.06304"
$ws.Range("E48").Value = "  -4.02%  "

# Row 49 - Decentraland
$ws.Range("D49").Value = "'0.4494"
$ws.Range("E49").Value = "  -5.61%  "

# Row 50 - NEARProtocol
$ws.Range("D50").Value = "'1.588"
$ws.Range("E50").Value = "  -4.42%  "

# Row 51 - Aave
$ws.Range("D51").Value = "'61.81"
$ws.Range("E51").Value = "  -4.53%  "
